$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is being restructured from a horizontal layout (row 1 =
# labels A1:S1, row 2 = values A2:S2, 19 columns) into a vertical
# layout (column A = labels, column B = values, 19 rows).

# Drop the old label/value columns C:S entirely - only A (labels) and
# B (values) survive into the new layout.
$ws.Range("C:S").EntireColumn.Delete()

# The old row-1 header cells (A1:S1) carried a bold/bordered/centered
# style. Propagate that same style down column A (the new labels
# column) for every label row, A1 already has it.
$ws.Range("A1").Copy()
$ws.Range("A2:A19").PasteSpecial(-4122)

# The old B1 header cell ("Price") had the header style too, but in
# the new layout B1 holds the plain value "Tesla, Inc." with no
# special formatting.
$ws.Range("B1").ClearFormats()

$labels = @(
    "Stock Name",
    "Price",
    "Change",
    "Previous Close",
    "Open",
    "Bid",
    "Ask",
    "Day's Range",
    "52 Week Range",
    "Volume",
    "Avg. Volume",
    "Market Cap",
    "Beta",
    "PE Ratio (TTM)",
    "EPS",
    "Earnings Date",
    "Forward Dividend & Yield",
    "Ex-Dividend Date",
    "1y Target Est"
)

$values = @(
    "Tesla, Inc.",
    262.67,
    -0.7799988,
    263.45,
    259.275,
    262.23,
    272.52,
    "250.73 - 266.2499",
    "138.8 - 488.54",
    101496505,
    86384041,
    844883361792,
    2.507,
    128.13171,
    2.05,
    "N/A",
    "N/A",
    "N/A",
    345.63025
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $labels[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
